{"js": "// The document has no word/styles.xml part yet (all paragraphs use the\n// implicit default \"Normal\" style). The edit adds an explicit style\n// definition for \"Normal\" to the document's style gallery, which causes\n// Word to materialize a word/styles.xml part (with the accompanying\n// content-type + relationship wiring) containing that single paragraph\n// style, without touching the body content itself.\nconst styles = context.document.getStyles();\nstyles.load(\"items\");\nawait context.sync();\n\nconst alreadyDefined = styles.items.some((s) => s.nameLocal === \"Normal\");\nif (!alreadyDefined) {\n  context.document.addStyle(\"Normal\", Word.StyleType.paragraph);\n  await context.sync();\n}\n", "ps1": "# The document has no word/styles.xml part yet (all paragraphs use the\n# implicit default \"Normal\" style). The edit adds an explicit style\n# definition for \"Normal\" to the document's style gallery, which causes\n# Word to materialize a word/styles.xml part (with the accompanying\n# content-type + relationship wiring) containing that single paragraph\n# style, without touching the body content itself.\n$d = $word.ActiveDocument\n\n$hasNormal = $false\nforeach ($s in $d.Styles) {\n    if ($s.NameLocal -eq \"Normal\") {\n        $hasNormal = $true\n        break\n    }\n}\n\nif (-not $hasNormal) {\n    $d.Styles.Add(\"Normal\", 1) | Out-Null\n}\n"}
